$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F8").Value = 1
$ws.Range("F10").Value = 6
$ws.Range("F14").Value = 5
$ws.Range("F16").Value = 0
$ws.Range("F17").Value = -1
$ws.Range("F19").Value = 1
$ws.Range("D20").Value = 12
$ws.Range("D24").Value = 17
$ws.Range("G24").Value = 6
$ws.Range("F25").Value = -1
$ws.Range("F35").Value = -1
$ws.Range("F36").Value = 2
$ws.Range("F49").Value = 0
$ws.Range("F53").Value = 1
$ws.Range("F56").Value = 3
$ws.Range("F57").Value = -1
$ws.Range("F58").Value = 1
$ws.Range("F62").Value = -2
$ws.Range("F71").Value = -4
